$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1163.3636
$ws.Range("I17").Value = 1034.3846
$ws.Range("J17").Value = 1349.6666
$ws.Range("K17").Value = 3103.1538
$ws.Range("L17").Value = 4048.9998
$ws.Range("M17").Value = -2935.1538
$ws.Range("N17").Value = -4384.9998
$ws.Range("H39").Value = 4321.7144
$ws.Range("I39").Value = 5020.3335
$ws.Range("J39").Value = 130
$ws.Range("K39").Value = 15061.0005
$ws.Range("L39").Value = 390
$ws.Range("M39").Value = -14765.0005
$ws.Range("N39").Value = -982
$ws.Range("H52").Value = 5000
$ws.Range("J52").Value = 5000
$ws.Range("L52").Value = 15000
$ws.Range("N52").Value = -15320
$ws.Range("H53").Value = 350
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H64").Value = 4560
$ws.Range("J64").Value = 4560
$ws.Range("L64").Value = 4560
$ws.Range("N64").Value = -5056
$ws.Range("H67").Value = 4560
$ws.Range("J67").Value = 4560
$ws.Range("L67").Value = 4560
$ws.Range("N67").Value = -6276
$ws.Range("H98").Value = 1612.7222
$ws.Range("I98").Value = 1612.7222
$ws.Range("K98").Value = 1612.7222
$ws.Range("M98").Value = -114.7221999999999
$ws.Range("H116").Value = 4999.5
$ws.Range("J116").Value = 4999.5
$ws.Range("L116").Value = 4999.5
$ws.Range("N116").Value = -11883.5
$ws.Range("H122").Value = 1612.7222
$ws.Range("I122").Value = 1612.7222
$ws.Range("K122").Value = 4838.1666
$ws.Range("M122").Value = -2388.1666
$ws.Range("H132").Value = 1166.6316
$ws.Range("I132").Value = 1170.4445
$ws.Range("K132").Value = 3511.3335
$ws.Range("M132").Value = -981.3335000000002
$ws.Range("H137").Value = 1869.1538
$ws.Range("J137").Value = 2214.2856
$ws.Range("L137").Value = 6642.8568
$ws.Range("N137").Value = -11742.8568
$ws.Range("H138").Value = 3410.8044
$ws.Range("I138").Value = 3216.077
$ws.Range("J138").Value = 3663.95
$ws.Range("K138").Value = 9648.231
$ws.Range("L138").Value = 10991.85
$ws.Range("M138").Value = -4508.231
$ws.Range("N138").Value = -21271.85

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 800
$ws.Range("I22").Value = 800
$ws.Range("K22").Value = 800
$ws.Range("M22").Value = -501
$ws.Range("H32").Value = 4426.0356
$ws.Range("I32").Value = 3025.8164
$ws.Range("K32").Value = 3025.8164
$ws.Range("M32").Value = -2738.8164
$ws.Range("H45").Value = 3068.75
$ws.Range("I45").Value = 2390.2
$ws.Range("K45").Value = 2390.2
$ws.Range("M45").Value = -2013.2
$ws.Range("H61").Value = 2084.9333
$ws.Range("I61").Value = 2043.3077
$ws.Range("K61").Value = 2043.3077
$ws.Range("M61").Value = -1831.3077
$ws.Range("H102").Value = 2370.5715
$ws.Range("I102").Value = 2098.8
$ws.Range("J102").Value = 3050
$ws.Range("K102").Value = 2098.8
$ws.Range("L102").Value = 3050
$ws.Range("M102").Value = -476.8000000000002
$ws.Range("N102").Value = -6294
$ws.Range("H132").Value = 1416
$ws.Range("I132").Value = 1468.3
$ws.Range("J132").Value = 893
$ws.Range("K132").Value = 4404.9
$ws.Range("L132").Value = 2679
$ws.Range("M132").Value = -1874.9
$ws.Range("N132").Value = -7739
$ws.Range("H136").Value = 2084.9333
$ws.Range("I136").Value = 2043.3077
$ws.Range("K136").Value = 6129.9231
$ws.Range("M136").Value = -3579.9231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3050.611
$ws.Range("I105").Value = 2549.75
$ws.Range("K105").Value = 2549.75
$ws.Range("M105").Value = -802.75
$ws.Range("H107").Value = 2013.1428
$ws.Range("I107").Value = 2348.75
$ws.Range("K107").Value = 2348.75
$ws.Range("M107").Value = -428.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3650.366
$ws.Range("I31").Value = 2269.1853
$ws.Range("J31").Value = 6314.0713
$ws.Range("K31").Value = 2269.1853
$ws.Range("L31").Value = 6314.0713
$ws.Range("M31").Value = -1974.1853
$ws.Range("N31").Value = -6904.0713
$ws.Range("H33").Value = 596
$ws.Range("I33").Value = 596
$ws.Range("K33").Value = 596
$ws.Range("M33").Value = -217
$ws.Range("H34").Value = 3650.366
$ws.Range("I34").Value = 2269.1853
$ws.Range("J34").Value = 6314.0713
$ws.Range("K34").Value = 2269.1853
$ws.Range("L34").Value = 6314.0713
$ws.Range("M34").Value = -2067.1853
$ws.Range("N34").Value = -6718.0713
$ws.Range("H37").Value = 15057
$ws.Range("J37").Value = 15057
$ws.Range("L37").Value = 15057
$ws.Range("N37").Value = -15271
$ws.Range("H105").Value = 3188.9167
$ws.Range("I105").Value = 853.6
$ws.Range("K105").Value = 853.6
$ws.Range("M105").Value = 893.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 5428.4287
$ws.Range("I137").Value = 7500
$ws.Range("J137").Value = 4599.8
$ws.Range("K137").Value = 22500
$ws.Range("L137").Value = 13799.4
$ws.Range("M137").Value = -17400
$ws.Range("N137").Value = -23999.4
$ws.Range("H139").Value = 1649.75
$ws.Range("I139").Value = 1649.75
$ws.Range("K139").Value = 4949.25
$ws.Range("M139").Value = 190.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 32176.572
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 32176.572
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 32176.572
$ws.Range("N26").Value = -32736.572
$ws.Range("M26").ClearContents()
$ws.Range("H35").Value = 7500000
$ws.Range("I35").Value = 10000000
$ws.Range("J35").Value = 5000000
$ws.Range("K35").Value = 10000000
$ws.Range("L35").Value = 5000000
$ws.Range("M35").Value = -9999702
$ws.Range("N35").Value = -5000596
$ws.Range("H40").Value = 19998
$ws.Range("J40").Value = 19998
$ws.Range("L40").Value = 19998
$ws.Range("N40").Value = -20300
$ws.Range("H50").Value = 32176.572
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 32176.572
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 32176.572
$ws.Range("N50").Value = -33172.572
$ws.Range("M50").ClearContents()
$ws.Range("H54").Value = 4124.75
$ws.Range("I54").Value = 1500
$ws.Range("J54").Value = 4999.6665
$ws.Range("K54").Value = 1500
$ws.Range("L54").Value = 4999.6665
$ws.Range("N54").Value = -5779.6665
$ws.Range("M54").Value = -1110

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 20037
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 20037
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 20037
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -21163
$ws.Range("H49").Value = 20037
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 20037
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 20037
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -20331
$ws.Range("H130").Value = 84625
$ws.Range("J130").Value = 84625
$ws.Range("L130").Value = 84625
$ws.Range("N130").Value = -94665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 25000
$ws.Range("I38").Value = 25000
$ws.Range("K38").Value = 25000
$ws.Range("M38").Value = -24527
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H45").Value = 20872
$ws.Range("I45").Value = 14991
$ws.Range("J45").Value = 23812.5
$ws.Range("K45").Value = 14991
$ws.Range("L45").Value = 23812.5
$ws.Range("M45").Value = -14500
$ws.Range("N45").Value = -24794.5
$ws.Range("H49").Value = 424999.25
$ws.Range("I49").Value = 424999.25
$ws.Range("K49").Value = 424999.25
$ws.Range("M49").Value = -424769.25
$ws.Range("H126").Value = 1727.2
$ws.Range("I126").Value = 1761.5
$ws.Range("J126").Value = 1590
$ws.Range("K126").Value = 5284.5
$ws.Range("L126").Value = 4770
$ws.Range("M126").Value = -2814.5
$ws.Range("N126").Value = -9710
$ws.Range("H132").Value = 3153.8235
$ws.Range("I132").Value = 2470.3845
$ws.Range("K132").Value = 7411.1535
$ws.Range("M132").Value = -4881.1535
